$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.43527889251709,-5.793748378753662,33.75442886352539),
    @(3.264275550842285,-8.339526176452637,33.21075820922852),
    @(-2.304214000701904,-44.59181213378906,10.13076019287109),
    @(-2.632879257202148,-47.35254669189453,16.20354652404785),
    @(-20.92312431335449,-10.33131408691406,-20.38119125366211),
    @(-20.04123497009277,-3.740296125411987,-22.7781810760498),
    @(-13.83809471130371,-1.457606315612793,10.62628936767578),
    @(-15.24413108825684,-5.868594169616699,10.52801132202148),
    @(-68.24240875244141,-37.54502487182617,-25.38205528259277),
    @(-69.99400329589844,-36.24868774414063,-23.62584495544434),
    @(3.574422836303711,-10.55799293518066,-1.943589687347412),
    @(1.990053176879883,-4.801568031311035,-0.5777735710144043),
    @(15.68469429016113,14.39326095581055,3.627434253692627),
    @(15.73915863037109,7.167607307434082,5.527215480804443),
    @(-0.1553750038146972,10.89957809448242,-15.06903266906738),
    @(-0.4681458473205566,11.8446626663208,-15.58957099914551),
    @(9.029394149780272,-76.68982696533203,2.083036422729492),
    @(9.611247062683104,-76.88915252685547,1.387777328491211),
    @(-11.94809913635254,10.77702140808106,4.135905742645264),
    @(-66.71106719970703,-81.49526214599609,59.02106475830078),
    @(-43.36343383789063,-8.654863357543945,-11.20723152160644),
    @(-40.74956512451172,-3.842054843902588,-15.99789905548096),
    @(-8.164802551269531,7.02423906326294,-4.845988273620605),
    @(-8.398514747619629,4.322819232940674,-5.664791107177734),
    @(-79.12488555908203,-28.08302688598633,-31.37481689453125),
    @(-81.6524658203125,-26.21007919311523,-27.36414909362793),
    @(-9.369277954101562,-24.22552871704102,-6.347414493560791),
    @(-8.832977294921875,-17.93884468078613,-7.180578708648682),
    @(-12.03362846374512,20.78299713134766,6.19299840927124),
    @(-11.60475063323975,14.68077850341797,6.681206226348877)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
